$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data source column (D) for rows 2-5 changes from "recovered" to "recovered_host"
$ws.Range("D2:D5").Value = "recovered_host"

# Update the saved selection to D5 (as in the edited file)
$ws.Range("D5").Select()
